$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Nominal Interest Rate" (A4) to "Nominal Interest Rate Reference"
$ws.Range("A4").Value = "Nominal Interest Rate Reference"

# Remove the "Volatility nominal IR" row (row 7) entirely, shifting rows below up
$ws.Rows.Item(7).Delete()

# Remove the trailing "Mean reversion nominal" row (now row 9 after the shift above)
$ws.Rows.Item(9).Delete()

# Widen column A to fit the new, longer label
$ws.Columns.Item(1).ColumnWidth = 28.7

# Update the active selection to A5, matching the saved view state
$ws.Range("A5").Select()
